# This workbook ("LinearVisualizer") drives a big table of formulas off of
# volatile functions (RAND() / RANDBETWEEN()), so simply recalculating the
# workbook re-rolls every dependent cell's cached value throughout the sheet.
# That recalculation cascade is what produced the bulk of the diff. On top of
# that, the sheet's window was scrolled/zoomed and a different cell was
# selected before the file was saved again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Force a full recalculation of every formula (including volatile ones),
# re-rolling RAND()/RANDBETWEEN() driven cells across the sheet.
$excel.CalculateFullRebuild()

# Update the view: zoom out to 25%, clear the scrolled-away top-left cell,
# and move the active selection to G108.
$excel.ActiveWindow.Zoom = 25
$ws.Range("G108").Select()
